# Update "想去人数" (F column) figures for the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1196
$wsExhibit.Range("F3").Value = 424
$wsExhibit.Range("F4").Value = 281
$wsExhibit.Range("F6").Value = 15
$wsExhibit.Range("F7").Value = 12348
$wsExhibit.Range("F9").Value = 19
$wsExhibit.Range("F11").Value = 161
$wsExhibit.Range("F12").Value = 12173
$wsExhibit.Range("F13").Value = 4835
$wsExhibit.Range("F14").Value = 4702
$wsExhibit.Range("F15").Value = 132
$wsExhibit.Range("F17").Value = 423
$wsExhibit.Range("F20").Value = 4
$wsExhibit.Range("F21").Value = 362
$wsExhibit.Range("F22").Value = 168

# Sheet "全部类型" (All types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1196
$wsAll.Range("F3").Value = 424
$wsAll.Range("F4").Value = 281
$wsAll.Range("F8").Value = 15
$wsAll.Range("F9").Value = 12348
$wsAll.Range("F11").Value = 19
$wsAll.Range("F13").Value = 161
$wsAll.Range("F14").Value = 12173
$wsAll.Range("F15").Value = 4835
$wsAll.Range("F16").Value = 4702
$wsAll.Range("F17").Value = 132
$wsAll.Range("F19").Value = 423
$wsAll.Range("F22").Value = 4
$wsAll.Range("F23").Value = 362
$wsAll.Range("F24").Value = 168
